$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "22"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44000.00"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "87"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "315698.00"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "455"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1173510.82"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "238"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "517041.00"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "824"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3141614.81"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "177"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "544816.18"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "105"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "252800.00"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "138"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "628217.26"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "326"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1366916.35"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "92"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "243987.00"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "505"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "2236814.03"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "9"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "33500.00"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "74"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "265072.36"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "75"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "342136.08"

$ws.Range("C122").NumberFormat = "@"
$ws.Range("C122").Value = "257"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("D122").Value = "717508.00"

$ws.Range("C123").NumberFormat = "@"
$ws.Range("C123").Value = "133"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "349012.45"

$ws.Range("C124").NumberFormat = "@"
$ws.Range("C124").Value = "520"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("D124").Value = "2376356.06"

$ws.Range("C196").NumberFormat = "@"
$ws.Range("C196").Value = "56"
$ws.Range("D196").NumberFormat = "@"
$ws.Range("D196").Value = "178300.00"

$ws.Range("C201").NumberFormat = "@"
$ws.Range("C201").Value = "670"
$ws.Range("D201").NumberFormat = "@"
$ws.Range("D201").Value = "2574882.58"

$ws.Range("C205").NumberFormat = "@"
$ws.Range("C205").Value = "157"
$ws.Range("D205").NumberFormat = "@"
$ws.Range("D205").Value = "491133.00"

$ws.Range("C209").NumberFormat = "@"
$ws.Range("C209").Value = "129"
$ws.Range("D209").NumberFormat = "@"
$ws.Range("D209").Value = "599185.14"

$ws.Range("C245").NumberFormat = "@"
$ws.Range("C245").Value = "118"
$ws.Range("D245").NumberFormat = "@"
$ws.Range("D245").Value = "344727.11"
